$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 1.111105389870137

# Delete rows 3 through 17 (old data no longer needed)
$ws.Range("A3:B17").EntireRow.Delete()
